$wb = $excel.ActiveWorkbook

$wsAdd = $wb.Worksheets.Item("AddProduct")
$wsGuest = $wb.Worksheets.Item("GuestUser")

# --- Content changes ---
# Rename header "PostalCode" -> "ZipCode" on GuestUser sheet
$wsGuest.Range("M1").Value = "ZipCode"

# Fix typo "VietName" -> "Viet Nam" on GuestUser sheet
$wsGuest.Range("O2").Value = "Viet Nam"

# Force the zip value in S2 to be stored as text (quote-prefixed), matching
# the target representation, while keeping the displayed value unchanged
$wsGuest.Range("S2").Value = "'83123456"

# --- Selection / active sheet changes ---
# Move selection on AddProduct sheet
$wsAdd.Range("D7").Select() | Out-Null

# Move selection on GuestUser sheet and activate it (becomes the active tab)
$wsGuest.Activate()
$wsGuest.Range("C1").Select() | Out-Null
